$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": rename project path sphr -> phr, and bump the Date ---
# (B2 and B17 happen to share the exact same underlying text/string, just like
#  in the source workbook, so both must be updated together.)
$metaWs = $wb.Worksheets.Item("Metadata")
$metaWs.Range("B2").Value = "http://hl7.org/fhir/uv/phr/StructureDefinition/AppleHealthKitHealthStoreCharacteristic"
$metaWs.Range("B17").Value = "http://hl7.org/fhir/uv/phr/StructureDefinition/AppleHealthKitHealthStoreCharacteristic"
$metaWs.Range("B8").Value = "2024-12-13T21:36:12-06:00"

# --- Sheet "Elements": rename project path sphr -> phr in the two Binding Value Set URLs ---
$elemWs = $wb.Worksheets.Item("Elements")
$elemWs.Range("Z3").Value = "http://hl7.org/fhir/uv/phr/ValueSet/apple-health-kit-characteristic-type-value-set"
$elemWs.Range("Z4").Value = "http://hl7.org/fhir/uv/phr/ValueSet/apple-health-kit-biological-sex-value-set"

# The shorter URL text narrows the best-fit width of column Z (26) from 77.21875 to
# 76.28515625 (OOXML "width" units). This COM surface only lets us drive that through
# the character-based ColumnWidth property (closest achievable value).
$elemWs.Columns.Item(26).ColumnWidth = 75.5
